{"js": "// Fix weird (doubled) spacing in the publications list on the resume:\n// \"...Hansman.  AIAA JSR, ...\" -> \"...Hansman. AIAA JSR, ...\"\n// i.e. collapse the duplicated space that preceded \"AIAA JSR\" in each\n// publication entry.\n\nconst body = context.document.body;\n\n// Locate every spot where two spaces are immediately followed by\n// \"AIAA JSR\" (there are two such publication entries on the resume).\nconst results = body.search(\"  AIAA JSR\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nconst fixedRanges = [];\nfor (let i = 0; i < results.items.length; i++) {\n  // Replace the double space + \"AIAA JSR\" with a single space + \"AIAA JSR\".\n  const newRange = results.items[i].insertText(\" AIAA JSR\", \"Replace\");\n  newRange.font.load(\"italic\");\n  fixedRanges.push(newRange);\n}\nawait context.sync();\n\n// The replacement can pick up stray italic formatting from the text that\n// immediately precedes the match (e.g. a preceding run ending in an italic\n// \".\"). \"AIAA JSR\" itself is never italic in the source, so normalize it\n// back if that happened.\nfor (const r of fixedRanges) {\n  if (r.font.italic) {\n    r.font.italic = false;\n  }\n}\nawait context.sync();\n", "ps1": "# Fix weird (doubled) spacing in the publications list on the resume:\n#   \"...Hansman.  AIAA JSR, ...\" -> \"...Hansman. AIAA JSR, ...\"\n# i.e. collapse the duplicated space that precedes \"AIAA JSR\" in each\n# publication entry.\n\n$d = $word.ActiveDocument\n\n$target = \"  AIAA JSR\"\n$replacement = \" AIAA JSR\"\n\n$searchStart = 0\n$docEnd = $d.Content.End\n$count = 0\n$maxIter = 50\n\nwhile ($count -lt $maxIter) {\n    $rng = $d.Range($searchStart, $docEnd)\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $target\n    $rng.Find.Forward = $true\n    $rng.Find.MatchCase = $true\n    $rng.Find.Wrap = 0   # wdFindStop - do not wrap back to the top\n\n    $found = $rng.Find.Execute()\n    if (-not $found) {\n        break\n    }\n\n    $count += 1\n\n    # Replace just the matched text (double space + \"AIAA JSR\") with a\n    # single space + \"AIAA JSR\".\n    $rng.Text = $replacement\n\n    # The in-place text replacement can pick up stray italic formatting\n    # from whatever immediately precedes the match (e.g. a preceding run\n    # ending in an italic \".\"). \"AIAA JSR\" itself is never italic in the\n    # source, so normalize it back if that happened.\n    if ($rng.Font.Italic) {\n        $rng.Font.Italic = $false\n    }\n\n    $searchStart = $rng.End\n    $docEnd = $d.Content.End\n}\n"}
